$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.004591068718582392
$ws.Cells.Item(2, 2).Value = 2.186912059783936
$ws.Cells.Item(2, 3).Value = 0.005427035968750715
$ws.Cells.Item(2, 4).Value = 6.246422290802002
$ws.Cells.Item(2, 5).Value = 4.384115219116211
$ws.Cells.Item(2, 6).Value = 3.762005805969238
$ws.Cells.Item(2, 7).Value = 1.004564762115479
$ws.Cells.Item(2, 8).Value = 1.314486384391785
$ws.Cells.Item(2, 9).Value = 1.938299298286438
$ws.Cells.Item(2, 10).Value = 3.32776141166687
$ws.Cells.Item(2, 11).Value = 5.660001277923584
$ws.Cells.Item(2, 12).Value = 5.990503787994385
$ws.Cells.Item(2, 13).Value = 3.653046369552612
$ws.Cells.Item(2, 14).Value = 0.004531692247837782
$ws.Cells.Item(2, 15).Value = 2.176706075668335
$ws.Cells.Item(2, 16).Value = 3.967877626419067
$ws.Cells.Item(2, 17).Value = 0.004104952327907085
$ws.Cells.Item(2, 18).Value = 2.203673601150513
$ws.Cells.Item(2, 19).Value = 0.1313426792621613
$ws.Cells.Item(2, 20).Value = 0.004068870097398758
$ws.Cells.Item(2, 21).Value = 0.00507756182923913
$ws.Cells.Item(2, 22).Value = 2.136415958404541
$ws.Cells.Item(2, 23).Value = 0.004902273416519165
$ws.Cells.Item(2, 24).Value = 0.004691822919994593
$ws.Cells.Item(2, 25).Value = 1.613394618034363
$ws.Cells.Item(3, 1).Value = 5.713008880615234
$ws.Cells.Item(3, 2).Value = 0.0002102738944813609
$ws.Cells.Item(3, 3).Value = 0.00244306237436831
$ws.Cells.Item(3, 4).Value = 0.0006176214665174484
$ws.Cells.Item(3, 5).Value = 1.905981540679932
$ws.Cells.Item(3, 6).Value = 0.0005943316500633955
$ws.Cells.Item(3, 7).Value = 0.001702116569504142
$ws.Cells.Item(3, 8).Value = 0.003019219264388084
$ws.Cells.Item(3, 9).Value = 1.835588574409485
$ws.Cells.Item(3, 10).Value = 0.02475676685571671
$ws.Cells.Item(3, 11).Value = 0.003329100087285042
$ws.Cells.Item(3, 12).Value = 1.283533811569214
$ws.Cells.Item(3, 13).Value = 0.7580472826957703
$ws.Cells.Item(3, 14).Value = 0.00204353523440659
$ws.Cells.Item(3, 15).Value = 0.187361404299736
$ws.Cells.Item(3, 16).Value = 0.005680469796061516
$ws.Cells.Item(3, 17).Value = 5.431578636169434
$ws.Cells.Item(3, 18).Value = 0.004768955521285534
$ws.Cells.Item(3, 19).Value = 0.6577330827713013
$ws.Cells.Item(3, 20).Value = 6.873856544494629
$ws.Cells.Item(3, 21).Value = 5.204483032226562
$ws.Cells.Item(3, 22).Value = 6.632120609283447
$ws.Cells.Item(3, 23).Value = 4.60166072845459
$ws.Cells.Item(3, 24).Value = 6.218491077423096
$ws.Cells.Item(3, 25).Value = 0.001210002461448312
$ws.Cells.Item(4, 1).Value = 1.9506756067276
$ws.Cells.Item(4, 2).Value = 0.004224009811878204
$ws.Cells.Item(4, 3).Value = 0.01214164029806852
$ws.Cells.Item(4, 4).Value = 0.01354831922799349
$ws.Cells.Item(4, 5).Value = 3.908854484558105
$ws.Cells.Item(4, 6).Value = 0.00608730036765337
$ws.Cells.Item(4, 7).Value = 2.115115165710449
$ws.Cells.Item(4, 8).Value = 8.032851219177246
$ws.Cells.Item(4, 9).Value = 4.617366790771484
$ws.Cells.Item(4, 10).Value = 1.840412259101868
$ws.Cells.Item(4, 11).Value = 0.003356289118528366
$ws.Cells.Item(4, 12).Value = 4.469799518585205
$ws.Cells.Item(4, 13).Value = 5.621973991394043
$ws.Cells.Item(4, 14).Value = 0.621628999710083
$ws.Cells.Item(4, 15).Value = 2.20595383644104
$ws.Cells.Item(4, 16).Value = 1.223170399665833
$ws.Cells.Item(4, 17).Value = 4.591441631317139
$ws.Cells.Item(4, 18).Value = 3.292394399642944
$ws.Cells.Item(4, 19).Value = 0.004453139379620552
$ws.Cells.Item(4, 20).Value = 0.009770682081580162
$ws.Cells.Item(4, 21).Value = 0.2070891261100769
$ws.Cells.Item(4, 22).Value = 0.684496283531189
$ws.Cells.Item(4, 23).Value = 5.382015228271484
$ws.Cells.Item(4, 24).Value = 0.8915483951568604
$ws.Cells.Item(4, 25).Value = 0.002153449924662709
$ws.Cells.Item(5, 1).Value = 1.461140275001526
$ws.Cells.Item(5, 2).Value = 0.004304415080696344
$ws.Cells.Item(5, 3).Value = 2.584566593170166
$ws.Cells.Item(5, 4).Value = 0.0201956033706665
$ws.Cells.Item(5, 5).Value = 3.477282285690308
$ws.Cells.Item(5, 6).Value = 0.003808443201705813
$ws.Cells.Item(5, 7).Value = 0.003639432601630688
$ws.Cells.Item(5, 8).Value = 6.121278762817383
$ws.Cells.Item(5, 9).Value = 0.003637688234448433
$ws.Cells.Item(5, 10).Value = 3.464127063751221
$ws.Cells.Item(5, 11).Value = 0.00824928842484951
$ws.Cells.Item(5, 12).Value = 0.003916628658771515
$ws.Cells.Item(5, 13).Value = 4.584489345550537
$ws.Cells.Item(5, 14).Value = 1.48308253288269
$ws.Cells.Item(5, 15).Value = 4.691103935241699
$ws.Cells.Item(5, 16).Value = 2.149705648422241
$ws.Cells.Item(5, 17).Value = 0.002888125367462635
$ws.Cells.Item(5, 18).Value = 4.577437877655029
$ws.Cells.Item(5, 19).Value = 5.123556137084961
$ws.Cells.Item(5, 20).Value = 2.871350526809692
$ws.Cells.Item(5, 21).Value = 3.953562259674072
$ws.Cells.Item(5, 22).Value = 4.243963718414307
$ws.Cells.Item(5, 23).Value = 0.003135476727038622
$ws.Cells.Item(5, 24).Value = 0.5569941997528076
$ws.Cells.Item(5, 25).Value = 0.00924244336783886
$ws.Cells.Item(6, 1).Value = 3.410658597946167
$ws.Cells.Item(6, 2).Value = 0.5193859934806824
$ws.Cells.Item(6, 3).Value = 3.805002689361572
$ws.Cells.Item(6, 4).Value = 6.188418388366699
$ws.Cells.Item(6, 5).Value = 0.03429876640439034
$ws.Cells.Item(6, 6).Value = 0.2813860177993774
$ws.Cells.Item(6, 7).Value = 5.831648349761963
$ws.Cells.Item(6, 8).Value = 0.5426952242851257
$ws.Cells.Item(6, 9).Value = 0.01319483108818531
$ws.Cells.Item(6, 10).Value = 1.831418633460999
$ws.Cells.Item(6, 11).Value = 5.579594612121582
$ws.Cells.Item(6, 12).Value = 0.00465824268758297
$ws.Cells.Item(6, 13).Value = 0.004029155243188143
$ws.Cells.Item(6, 14).Value = 5.427706241607666
$ws.Cells.Item(6, 15).Value = 1.787075757980347
$ws.Cells.Item(6, 16).Value = 0.003378417575731874
$ws.Cells.Item(6, 17).Value = 5.225881576538086
$ws.Cells.Item(6, 18).Value = 0.007809703703969717
$ws.Cells.Item(6, 19).Value = 0.5318001508712769
$ws.Cells.Item(6, 20).Value = 0.01222683489322662
$ws.Cells.Item(6, 21).Value = 1.937697172164917
$ws.Cells.Item(6, 22).Value = 0.004065602086484432
$ws.Cells.Item(6, 23).Value = 1.937944531440735
$ws.Cells.Item(6, 24).Value = 0.0103471614420414
$ws.Cells.Item(6, 25).Value = 5.803413391113281
$ws.Cells.Item(7, 1).Value = 2.844599008560181
$ws.Cells.Item(7, 2).Value = 4.741164207458496
$ws.Cells.Item(7, 3).Value = 0.002243695314973593
$ws.Cells.Item(7, 4).Value = 1.268503427505493
$ws.Cells.Item(7, 5).Value = 8.398907661437988
$ws.Cells.Item(7, 6).Value = 4.01053524017334
$ws.Cells.Item(7, 7).Value = 0.005037608556449413
$ws.Cells.Item(7, 8).Value = 2.781199932098389
$ws.Cells.Item(7, 9).Value = 0.001711501041427255
$ws.Cells.Item(7, 10).Value = 0.001031458843499422
$ws.Cells.Item(7, 11).Value = 2.740577936172485
$ws.Cells.Item(7, 12).Value = 1.787148714065552
$ws.Cells.Item(7, 13).Value = 0.4003611207008362
$ws.Cells.Item(7, 14).Value = 0.005779871717095375
$ws.Cells.Item(7, 15).Value = 0.001129458658397198
$ws.Cells.Item(7, 16).Value = 1.951138734817505
$ws.Cells.Item(7, 17).Value = 0.03499886021018028
$ws.Cells.Item(7, 18).Value = 1.159641742706299
$ws.Cells.Item(7, 19).Value = 6.422484397888184
$ws.Cells.Item(7, 20).Value = 0.007131356745958328
$ws.Cells.Item(7, 21).Value = 5.555715084075928
$ws.Cells.Item(7, 22).Value = 1.059916973114014
$ws.Cells.Item(7, 23).Value = 0.130394920706749
$ws.Cells.Item(7, 24).Value = 0.008469589985907078
$ws.Cells.Item(7, 25).Value = 5.590366363525391
$ws.Cells.Item(8, 1).Value = 1.498386859893799
$ws.Cells.Item(8, 2).Value = 0.233984500169754
$ws.Cells.Item(8, 3).Value = 0.3373751044273376
$ws.Cells.Item(8, 4).Value = 4.206531524658203
$ws.Cells.Item(8, 5).Value = 0.0001562673423904926
$ws.Cells.Item(8, 6).Value = 0.01127286721020937
$ws.Cells.Item(8, 7).Value = 1.683412313461304
$ws.Cells.Item(8, 8).Value = 0.001065192511305213
$ws.Cells.Item(8, 9).Value = 5.496126651763916
$ws.Cells.Item(8, 10).Value = 0.001083079841919243
$ws.Cells.Item(8, 11).Value = 5.948645114898682
$ws.Cells.Item(8, 12).Value = 5.852470874786377
$ws.Cells.Item(8, 13).Value = 0.001966696232557297
$ws.Cells.Item(8, 14).Value = 0.2667432427406311
$ws.Cells.Item(8, 15).Value = 0.0009038528660312295
$ws.Cells.Item(8, 16).Value = 0.4684150516986847
$ws.Cells.Item(8, 17).Value = 2.199663400650024
$ws.Cells.Item(8, 18).Value = 8.177361488342285
$ws.Cells.Item(8, 19).Value = 0.07305580377578735
$ws.Cells.Item(8, 20).Value = 0.1000940054655075
$ws.Cells.Item(8, 21).Value = 0.886012852191925
$ws.Cells.Item(8, 22).Value = 0.0004968990106135607
$ws.Cells.Item(8, 23).Value = 5.005722999572754
$ws.Cells.Item(8, 24).Value = 5.099513530731201
$ws.Cells.Item(8, 25).Value = 4.210269927978516
$ws.Cells.Item(9, 1).Value = 0.005901766940951347
$ws.Cells.Item(9, 2).Value = 4.256443500518799
$ws.Cells.Item(9, 3).Value = 5.195938110351562
$ws.Cells.Item(9, 4).Value = 0.006182674318552017
$ws.Cells.Item(9, 5).Value = 0.001698993495665491
$ws.Cells.Item(9, 6).Value = 3.949512004852295
$ws.Cells.Item(9, 7).Value = 0.3545295894145966
$ws.Cells.Item(9, 8).Value = 0.8364130854606628
$ws.Cells.Item(9, 9).Value = 0.09276384860277176
$ws.Cells.Item(9, 10).Value = 6.77872896194458
$ws.Cells.Item(9, 11).Value = 0.01233084686100483
$ws.Cells.Item(9, 12).Value = 0.004123758990317583
$ws.Cells.Item(9, 13).Value = 0.003267535008490086
$ws.Cells.Item(9, 14).Value = 5.680654525756836
$ws.Cells.Item(9, 15).Value = 7.149709224700928
$ws.Cells.Item(9, 16).Value = 4.894673347473145
$ws.Cells.Item(9, 17).Value = 0.06682708114385605
$ws.Cells.Item(9, 18).Value = 0.007031361106783152
$ws.Cells.Item(9, 19).Value = 0.008566330187022686
$ws.Cells.Item(9, 20).Value = 0.3613593876361847
$ws.Cells.Item(9, 21).Value = 0.004120903089642525
$ws.Cells.Item(9, 22).Value = 7.000304222106934
$ws.Cells.Item(9, 23).Value = 2.666243314743042
$ws.Cells.Item(9, 24).Value = 0.01006762124598026
$ws.Cells.Item(9, 25).Value = 0.003580531338229775
$ws.Cells.Item(10, 1).Value = 3.255012989044189
$ws.Cells.Item(10, 2).Value = 5.15381908416748
$ws.Cells.Item(10, 3).Value = 2.272378206253052
$ws.Cells.Item(10, 4).Value = 0.003712342819198966
$ws.Cells.Item(10, 5).Value = 2.640583992004395
$ws.Cells.Item(10, 6).Value = 4.243880748748779
$ws.Cells.Item(10, 7).Value = 0.004041530191898346
$ws.Cells.Item(10, 8).Value = 0.007993529550731182
$ws.Cells.Item(10, 9).Value = 4.409947872161865
$ws.Cells.Item(10, 10).Value = 2.664298295974731
$ws.Cells.Item(10, 11).Value = 0.004182351287454367
$ws.Cells.Item(10, 12).Value = 4.112948417663574
$ws.Cells.Item(10, 13).Value = 5.114810943603516
$ws.Cells.Item(10, 14).Value = 0.004202362149953842
$ws.Cells.Item(10, 15).Value = 1.058089256286621
$ws.Cells.Item(10, 16).Value = 0.4905550181865692
$ws.Cells.Item(10, 17).Value = 3.316354751586914
$ws.Cells.Item(10, 18).Value = 0.002530722878873348
$ws.Cells.Item(10, 19).Value = 6.137276649475098
$ws.Cells.Item(10, 20).Value = 0.643217146396637
$ws.Cells.Item(10, 21).Value = 2.728050231933594
$ws.Cells.Item(10, 22).Value = 0.002977384487167001
$ws.Cells.Item(10, 23).Value = 1.678724646568298
$ws.Cells.Item(10, 24).Value = 1.100385785102844
$ws.Cells.Item(10, 25).Value = 1.058057188987732
$ws.Cells.Item(11, 1).Value = 4.212934494018555
$ws.Cells.Item(11, 2).Value = 2.952179431915283
$ws.Cells.Item(11, 3).Value = 5.61153507232666
$ws.Cells.Item(11, 4).Value = 3.482506990432739
$ws.Cells.Item(11, 5).Value = 0.0003373771905899048
$ws.Cells.Item(11, 6).Value = 3.166513681411743
$ws.Cells.Item(11, 7).Value = 3.423604488372803
$ws.Cells.Item(11, 8).Value = 0.0001673406222835183
$ws.Cells.Item(11, 9).Value = 0.04144493117928505
$ws.Cells.Item(11, 10).Value = 2.198819160461426
$ws.Cells.Item(11, 11).Value = 0.00232517090626061
$ws.Cells.Item(11, 12).Value = 0.0002239847090095282
$ws.Cells.Item(11, 13).Value = 0.001565230544656515
$ws.Cells.Item(11, 14).Value = 4.012048721313477
$ws.Cells.Item(11, 15).Value = 2.243364334106445
$ws.Cells.Item(11, 16).Value = 8.73178768157959
$ws.Cells.Item(11, 17).Value = 2.383157253265381
$ws.Cells.Item(11, 18).Value = 0.0007974395994096994
$ws.Cells.Item(11, 19).Value = 0.0002751208376139402
$ws.Cells.Item(11, 20).Value = 1.319719552993774
$ws.Cells.Item(11, 21).Value = 0.9394583106040955
$ws.Cells.Item(11, 22).Value = 0.003731885924935341
$ws.Cells.Item(11, 23).Value = 0.01734456792473793
$ws.Cells.Item(11, 24).Value = 0.004488531500101089
$ws.Cells.Item(11, 25).Value = 5.287984848022461
